$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("2023-10-02 16:23:58", "hatespeech", "def", 10, 0.3550838294873909),
    @("2023-10-02 16:24:05", "hatespeech", "def", 20, 0.2604653727803925),
    @("2023-10-02 16:24:05", "hatespeech", "def", 30, 0.2036138335506539),
    @("2023-10-02 16:24:05", "hatespeech", "def", 40, 0.1612300686123467)
)

$startRow = 28
for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $startRow + $i
    $ws.Cells.Item($row, 1).Value = $data[$i][0]
    $ws.Cells.Item($row, 2).Value = $data[$i][1]
    $ws.Cells.Item($row, 3).Value = $data[$i][2]
    $ws.Cells.Item($row, 4).Value = $data[$i][3]
    $ws.Cells.Item($row, 5).Value = $data[$i][4]
}
